$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string for row 7, column C
$newText = "Nová šablona - přesunutí JS, index - menu, patička, mapa, nový GIT repozitář"

# Fill in the new row of data (row 7)
$ws.Range("A7").Value = 42885
$ws.Range("B7").Value = 2.5
$ws.Range("C7").Value = $newText

# Copy date formatting from A6 down to A7 so it keeps the same style (s="1")
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 42885

# Apply a two-decimal number format to column B data cells (B2:B7),
# matching the centered alignment already used there.
$ws.Range("B2:B7").NumberFormat = "0.00"
$ws.Range("B2:B7").HorizontalAlignment = -4108

# Update the active selection to the new last cell C7
$ws.Range("C7").Select()

$excel.CalculateFull()
